$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("openbis-metadata")
$ws2 = $wb.Worksheets.Item("openbis-data")

# --- openbis-metadata (sheet1) ---
# B6/C6 text updates: units now include the new "fmol/ug protein digest" unit
$ws1.Range("B6").Value = "fmol/ug protein digest"
$ws1.Range("C6").Value = "One of mM, uM, Percent, RatioT1, RatioCs, or AU, Dimensionless, fmol/ug protein digest"

# Column C is now wider to fit the longer description text (bestFit/AutoFit applied)
$ws1.Columns.Item(3).ColumnWidth = 89.28571428571429

# Selection moved to C6
$ws1.Range("C6").Select()

# --- openbis-data (sheet2) ---
# Header row now encodes the Bio-replicate / Timepoint axes in the column name
$ws2.Range("C1").Value = "0::Mean::B1_B2::T1_T2"
$ws2.Range("D1").Value = "0::Std::B1_B2::T1_T2"
# E1/F1 begin with "+" so they are entered with a leading apostrophe (quote-prefixed text)
$ws2.Range("E1").Formula = "'+2100::Mean::B1_B2::T1_T2"
$ws2.Range("F1").Formula = "'+2100::Std::B1_B2::T1_T2"

# New columns C:F get explicit widths (bestFit on C, D, F)
$ws2.Columns.Item(3).ColumnWidth = 23.857142857142858
$ws2.Columns.Item(4).ColumnWidth = 16.857142857142858
$ws2.Columns.Item(5).ColumnWidth = 24.0
$ws2.Columns.Item(6).ColumnWidth = 20.857142857142858

# Selection moved to E6
$ws2.Range("E6").Select()
